# Update scenario number test case file
# Renumbers the "Scenarios" column (column B) on the TestCases sheet for
# rows 11-113 so the numbering is sequential / corrected, and moves the
# frozen-pane view / active selection down to the bottom of the updated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")
$ws.Activate()

$scenarioNumbers = @{
    11 = 6
    12 = 7
    13 = 8
    14 = 9
    15 = 10
    16 = 11
    17 = 12
    18 = 13
    19 = 15
    20 = 14
    21 = 16
    22 = 17
    23 = 18
    24 = 19
    25 = 20
    26 = 21
    27 = 22
    28 = 23
    29 = 24
    30 = 25
    31 = 26
    32 = 27
    33 = 28
    34 = 29
    35 = 30
    36 = 31
    37 = 32
    38 = 33
    39 = 34
    40 = 35
    41 = 36
    42 = 37
    43 = 38
    44 = 39
    45 = 40
    46 = 41
    47 = 42
    48 = 42
    49 = 43
    50 = 44
    51 = 45
    52 = 46
    53 = 47
    54 = 48
    55 = 49
    56 = 50
    57 = 51
    58 = 104
    59 = 52
    60 = 52
    61 = 53
    62 = 54
    63 = 55
    64 = 56
    65 = 58
    66 = 57
    67 = 59
    68 = 60
    69 = 60
    70 = 61
    71 = 62
    72 = 63
    73 = 64
    74 = 64
    75 = 65
    76 = 66
    77 = 67
    78 = 68
    79 = 69
    80 = 70
    81 = 71
    82 = 72
    83 = 73
    84 = 74
    85 = 75
    86 = 76
    87 = 77
    88 = 78
    89 = 79
    90 = 80
    91 = 81
    92 = 82
    93 = 83
    94 = 84
    95 = 85
    96 = 86
    97 = 87
    98 = 88
    99 = 89
    100 = 90
    101 = 91
    102 = 92
    103 = 93
    104 = 94
    105 = 95
    106 = 96
    107 = 97
    108 = 98
    109 = 99
    110 = 100
    111 = 101
    112 = 102
    113 = 103
}

foreach ($row in $scenarioNumbers.Keys) {
    $ws.Cells.Item($row, 2).Value = $scenarioNumbers[$row]
}

# Scroll the frozen view down to the bottom of the sheet and move the
# active selection to match where editing left off.
$ws.Range("A103").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("B118").Select()
